# Apply cryptos.xlsx price/volume updates (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.835.80"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "1.905.42"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "'313.00"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'0.5037"
$ws.Range("E7").Value = "  +4.72%  "
$ws.Range("D8").Value = "'0.3808"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "'0.07276"
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").Value = "'0.9060"
$ws.Range("E10").Value = "  -2.77%  "
$ws.Range("D11").Value = "'20.89"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").Value = "'0.07648"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "1.883.68"
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("D14").Value = "'5.488"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "'91.80"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "'0.000008716"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "27.864.19"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "'14.60"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").Value = "'6.579"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").Value = "'153.58"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").Value = "'1.879"
$ws.Range("E25").Value = "  -1.46%  "
$ws.Range("D26").Value = "'2.223"
$ws.Range("E26").Value = "  +4.41%  "
$ws.Range("D27").Value = "'18.36"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").Value = "'115.35"
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("D29").Value = "'4.902"
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("D30").Value = "'0.09010"
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").Value = "'3.210"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("E32").Value = "  -2.80%  "
$ws.Range("D33").Value = "'4.693"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").Value = "'0.7654"
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("D35").Value = "'0.02065"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").Value = "'2.519"
$ws.Range("E36").Value = "  -3.55%  "
$ws.Range("E37").Value = "  -1.81%  "
$ws.Range("D38").Value = "'0.5519"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("E39").Value = "  +0.60%  "
$ws.Range("D40").Value = "'0.05254"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").Value = "'6.873"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("D42").Value = "'8.437"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("D43").Value = "'0.1509"
$ws.Range("E43").Value = "  -1.43%  "
$ws.Range("D44").Value = "'111.22"
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("D45").Value = "'10.56"
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").Value = "'0.4798"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'1.626"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").Value = "'67.42"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "'0.06062"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("D51").Value = "'0.9023"
$ws.Range("E51").Value = "  +0.23%  "
